$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra rows (rows 6-16) that are no longer needed, keeping rows 1-5
$ws.Rows("6:16").Delete()

# Set the new header for column C
$ws.Range("C1").Value = "Service"

# Build the data-row template style on A2 (already has the thin border):
# add left/center alignment
$dataTemplate = $ws.Range("A2")
$dataTemplate.HorizontalAlignment = -4131
$dataTemplate.VerticalAlignment = -4108

# Propagate the data-row format (border + alignment) onto the whole data block,
# including the new column C
$dataTemplate.Copy()
$ws.Range("A2:C5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Build the header-row template style on A1 (already has the thin border):
# add left/center alignment plus a 50% gray themed fill
$headerTemplate = $ws.Range("A1")
$headerTemplate.HorizontalAlignment = -4131
$headerTemplate.VerticalAlignment = -4108
$headerTemplate.Interior.ThemeColor = 1
$headerTemplate.Interior.TintAndShade = 0.499984740745262

# Propagate the header format onto the whole header row, including column C
$headerTemplate.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Update the selected cell to match the edited workbook
$ws.Range("C8").Select()
